# Fill in the two new contact rows (14 and 15) that were added to the
# contacts sheet, and leave the selection on the last cell typed (E15) -
# matching the post-edit cursor position recorded in the saved workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14: new contact "Brian Lanning"
$ws.Range("C14").Value = "Brian "
$ws.Range("D14").Value = "Lanning"
$ws.Range("E14").Value = "8888888888"

# Row 15: new contact "another person"
$ws.Range("C15").Value = "another"
$ws.Range("D15").Value = "person"
$ws.Range("E15").Value = "988-708-9782"

# Leave the cursor on the last entered cell, scrolled back to the top of
# the sheet (no more topLeftCell override).
$ws.Range("E15").Select() | Out-Null
